# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-11-09 16:15:32
#
# The "Recorded By" column (column G) on the "Session Analysis Results"
# sheet lists the users who recorded / touched a session, separated by
# ", ". Upstream re-ordered that list (reversed it) for every row that
# has more than one recorder. Single-recorder rows are left untouched
# (reversing a one-element list is a no-op anyway).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$recordedByCol = 7   # column G ("Recorded By")

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    $val = $cell.Value2

    if ($val -ne $null -and $val -like "*,*") {
        $parts = $val -split ','
        for ($i = 0; $i -lt $parts.Length; $i++) {
            $parts[$i] = $parts[$i].Trim()
        }

        # Reverse the order of the recorder list.
        $count = $parts.Length
        $reversedParts = @()
        for ($i = $count - 1; $i -ge 0; $i--) {
            $reversedParts += $parts[$i]
        }

        $newVal = [string]::Join(', ', $reversedParts)
        $cell.Value = $newVal
    }
}
